# "Added Heroku in skills": split the trailing ", Tableau." run into
# ", Tableau" (unchanged run) + a new ", Heroku." run, Arial-formatted
# like its neighbours.
$d = $word.ActiveDocument
$r = $d.Content

$found = $r.Find.Execute(", Tableau.", $true, $false, $false, $false, $false, $true, 1, $false, ", Tableau", 2)

$r.Collapse(0)
$r.InsertAfter(", Heroku.")
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
